# Convert the query to its representation of the reduced dimension.
# The tfidf sheet gains a 5th document ("query.txt") which changes the
# document-frequency / idf weighting for every term, so every existing
# tf-idf value in the table is recomputed and a new row is appended for
# the query document.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reset the sheet view: unhide the helper header row/column and drop
#     the scrolled/selected state left over from editing ---
$ws.Rows.Item(1).Hidden = $false
$ws.Columns.Item(1).Hidden = $false
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A1").Select() | Out-Null

# --- Update the recomputed tf-idf values for the existing documents ---

# doc_1.txt (row 2)
$ws.Range("C2").Value = 0.1021651247531982
$ws.Range("D2").Value = 0.3218875824868201
$ws.Range("E2").Value = 0.3218875824868201
$ws.Range("F2").Value = 0.04462871026284196
$ws.Range("I2").Value = 0.183258146374831

# doc_2.txt (row 3)
$ws.Range("B3").Value = 0.4023594781085251
$ws.Range("C3").Value = 0.1277064059414977
$ws.Range("G3").Value = 0.1277064059414977
$ws.Range("J3").Value = 0.1277064059414977

# doc_3.txt (row 4)
$ws.Range("C4").Value = 0.1021651247531982
$ws.Range("F4").Value = 0.08925742052568392
$ws.Range("G4").Value = 0.1021651247531982
$ws.Range("J4").Value = 0.1021651247531982

# doc_4.txt (row 5)
$ws.Range("F5").Value = 0.04462871026284196
$ws.Range("H5").Value = 0.3218875824868201
$ws.Range("I5").Value = 0.3665162927496621
$ws.Range("K5").Value = 0.3218875824868201

# --- Append the new query.txt row (row 6) with its tf-idf weights ---

# Copy the formatting of the previous row label cell (bold/border/centered
# style) down onto the new label cell, then set its text.
$ws.Range("A5").Copy($ws.Range("A6"))
$ws.Range("A6").Value = "query.txt"

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0.07438118377140325
$ws.Range("G6").Value = 0.1702752079219969
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0.1702752079219969
$ws.Range("K6").Value = 0
